# Move the "Program-Year" column (AB) to the front of the table (new column A),
# shifting CATEGORY..Other Remarks one column to the right, and rename the
# REGION header (now in column C) to "Region".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("AB").Cut() | Out-Null
$ws.Columns("A").Insert() | Out-Null

$ws.Range("C1").Value = "Region"
